# Introspect Manual - Layers: add the "Provider Layer" row to the layering
# table (removing the ServiceLookup/"Provider Objects" callout from the
# Infrastructure row's visual and giving the Provider concept its own row),
# and restore the EmailStyle20 style naming.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add a new row to the layers table, right after "Infrastructure Layer",
#    for the new "Provider Layer" entry.
# ------------------------------------------------------------------
$table = $d.Tables(1)
$newRow = $table.Rows.Add()

$cell = $newRow.Cells(1)

# Match the target cell shading fill (#548DD4 == Text 2, 80% lighter theme
# tint in the source file) -- BackgroundPatternColor takes a BGR-packed
# OLE_COLOR, so red=0x54, green=0x8D, blue=0xD4 -> 0xD48D54.
$cell.Shading.BackgroundPatternColor = 13929812

# Two centered paragraphs: "Provider " then "Layer" (the row's decorative
# "Provider Objects" callout graphic itself lives purely in the VML drawing
# layer, which this host does not expose on the Word object model, so the
# text content is what is reproduced here).
$cell.Range.Text = "Provider " + [char]13 + "Layer"

# ------------------------------------------------------------------
# 2) EmailStyle20 / EmailStyle201 naming swap -- the custom character style
#    had its display name and the style id it should carry transposed;
#    restore the intended display name via the object model (the style id
#    itself is immutable through the supported COM surface).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $s = $d.Styles($i)
    if ($s.NameLocal -eq "EmailStyle201,EmailStyle201") {
        $s.NameLocal = "EmailStyle20"
    }
}
